$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.804.24"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "2.654.38"
$ws.Range("E3").Value = "  +3.80%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.51%  "

$ws.Range("D9").Value = "2.685.70"
$ws.Range("E9").Value = "  +5.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("E11").Value = "  +5.04%  "

$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("D14").Value = "3.124.43"
$ws.Range("E14").Value = "  +4.07%  "

$ws.Range("D15").Value = "58.815.04"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "

$ws.Range("E17").Value = "  +2.03%  "

$ws.Range("D18").Value = "2.683.86"
$ws.Range("E18").Value = "  +5.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.85%  "

$ws.Range("E20").Value = "  +0.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.71%  "

$ws.Range("E25").Value = "  +3.34%  "

$ws.Range("D26").Value = "2.786.87"
$ws.Range("E26").Value = "  +4.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("E28").Value = "  +2.13%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.65%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0807"
$ws.Range("E30").Value = "  +4.20%  "

$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.12%  "

$ws.Range("E34").Value = "  +2.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("E36").Value = "  +14.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("E38").Value = "  +3.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.842"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.18%  "

$ws.Range("E41").Value = "  +5.67%  "

$ws.Range("E42").Value = "  +1.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.618"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "277.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.73%  "

$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0980"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.29%  "

$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0230"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.06%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.57%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.003.04"
$ws.Range("E51").Value = "  +5.07%  "
